$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers
$ws.Range("A1").Value = "Índice"
$ws.Range("B1").Value = "Distancia"
$ws.Range("C1").Value = "max"
$ws.Range("D1").Value = "min"
$ws.Range("E1").Value = "Tempo"

# Data rows
$data = @(
    @(0, 5942.633333333333, 6432, 5483, 0.08182970682779948),
    @(1, 6612.933333333333, 7050, 5826, 0.08760472933451334),
    @(2, 6618.233333333334, 7019, 5873, 0.08665369351704916),
    @(3, 7255.133333333333, 7824, 6413, 0.08531359831492107),
    @(4, 5947.333333333333, 6433, 5259, 0.08318480650583902),
    @(5, 6521.066666666667, 7256, 5589, 0.0885109821955363),
    @(6, 6242.1, 6590, 5789, 0.08738417625427246),
    @(7, 6362.7, 6822, 5738, 0.08865559101104736),
    @(8, 5850.9, 6329, 5063, 0.08399654229482015),
    @(9, 6104.6, 6749, 5124, 0.08320171038309733)
)

$rowIndex = 2
foreach ($row in $data) {
    $ws.Cells.Item($rowIndex, 1).Value = $row[0]
    $ws.Cells.Item($rowIndex, 2).Value = $row[1]
    $ws.Cells.Item($rowIndex, 3).Value = $row[2]
    $ws.Cells.Item($rowIndex, 4).Value = $row[3]
    $ws.Cells.Item($rowIndex, 5).Value = $row[4]
    $rowIndex++
}
